$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3616.7727
$ws.Range("I98").Value = 1920.5
$ws.Range("J98").Value = 11250
$ws.Range("K98").Value = 1920.5
$ws.Range("L98").Value = 11250
$ws.Range("M98").Value = -422.5
$ws.Range("N98").Value = -14246

$ws.Range("H122").Value = 3616.7727
$ws.Range("I122").Value = 1920.5
$ws.Range("J122").Value = 11250
$ws.Range("K122").Value = 5761.5
$ws.Range("L122").Value = 33750
$ws.Range("M122").Value = -3311.5
$ws.Range("N122").Value = -38650

$ws.Range("H131").Value = 8368.333000000001
$ws.Range("I131").Value = 7500
$ws.Range("J131").Value = 10105
$ws.Range("K131").Value = 22500
$ws.Range("L131").Value = 30315
$ws.Range("M131").Value = -17460
$ws.Range("N131").Value = -40395

$ws.Range("H137").Value = 3352.4644
$ws.Range("I137").Value = 3143.5278
$ws.Range("J137").Value = 3728.55
$ws.Range("K137").Value = 9430.5834
$ws.Range("L137").Value = 11185.65
$ws.Range("M137").Value = -6880.5834
$ws.Range("N137").Value = -16285.65

$ws.Range("H138").Value = 3924.919
$ws.Range("I138").Value = 1656
$ws.Range("J138").Value = 4501.7627
$ws.Range("K138").Value = 4968
$ws.Range("L138").Value = 13505.2881
$ws.Range("M138").Value = 172
$ws.Range("N138").Value = -23785.2881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1598.8334
$ws.Range("I2").Value = 1793.3334
$ws.Range("J2").Value = 1404.3334
$ws.Range("K2").Value = 1793.3334
$ws.Range("L2").Value = 1404.3334
$ws.Range("M2").Value = -1680.3334
$ws.Range("N2").Value = -1630.3334

$ws.Range("H61").Value = 2523.889
$ws.Range("I61").Value = 1792.8889
$ws.Range("J61").Value = 3254.889
$ws.Range("K61").Value = 1792.8889
$ws.Range("L61").Value = 3254.889
$ws.Range("M61").Value = -1580.8889
$ws.Range("N61").Value = -3678.889

$ws.Range("H74").Value = 3899.4211
$ws.Range("I74").Value = 4313.6294
$ws.Range("K74").Value = 4313.6294
$ws.Range("M74").Value = -3439.6294

$ws.Range("H77").Value = 3899.4211
$ws.Range("I77").Value = 4313.6294
$ws.Range("K77").Value = 21568.147
$ws.Range("M77").Value = -17200.147

$ws.Range("H116").Value = 1598.8334
$ws.Range("I116").Value = 1793.3334
$ws.Range("J116").Value = 1404.3334
$ws.Range("K116").Value = 1793.3334
$ws.Range("L116").Value = 1404.3334
$ws.Range("M116").Value = 500.6666
$ws.Range("N116").Value = -5992.3334

$ws.Range("H122").Value = 4977.857
$ws.Range("I122").Value = 2006
$ws.Range("J122").Value = 6166.6
$ws.Range("K122").Value = 6018
$ws.Range("L122").Value = 18499.8
$ws.Range("M122").Value = -3568
$ws.Range("N122").Value = -23399.8

$ws.Range("H132").Value = 3231.3333
$ws.Range("I132").Value = 1706.7693
$ws.Range("J132").Value = 4222.3
$ws.Range("K132").Value = 5120.3079
$ws.Range("L132").Value = 12666.9
$ws.Range("M132").Value = -2590.3079
$ws.Range("N132").Value = -17726.9

$ws.Range("H136").Value = 2523.889
$ws.Range("I136").Value = 1792.8889
$ws.Range("J136").Value = 3254.889
$ws.Range("K136").Value = 5378.6667
$ws.Range("L136").Value = 9764.667000000001
$ws.Range("M136").Value = -2828.6667
$ws.Range("N136").Value = -14864.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1598.8334
$ws.Range("I3").Value = 1793.3334
$ws.Range("J3").Value = 1404.3334
$ws.Range("K3").Value = 1793.3334
$ws.Range("L3").Value = 1404.3334
$ws.Range("M3").Value = -1679.3334
$ws.Range("N3").Value = -1632.3334

$ws.Range("H43").Value = 79800
$ws.Range("J43").Value = 79800
$ws.Range("L43").Value = 79800
$ws.Range("N43").Value = -80162

$ws.Range("H105").Value = 2155.0454
$ws.Range("I105").Value = 1970.5883
$ws.Range("J105").Value = 2782.2
$ws.Range("K105").Value = 1970.5883
$ws.Range("L105").Value = 2782.2
$ws.Range("M105").Value = -223.5882999999999
$ws.Range("N105").Value = -6276.2

$ws.Range("H134").Value = 2852.25
$ws.Range("I134").Value = 1578
$ws.Range("J134").Value = 5995.4
$ws.Range("K134").Value = 4734
$ws.Range("L134").Value = 17986.2
$ws.Range("M134").Value = -2199
$ws.Range("N134").Value = -23056.2

$ws.Range("H137").Value = 50778
$ws.Range("J137").Value = 50778
$ws.Range("L137").Value = 50778
$ws.Range("N137").Value = -60978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4230.7666
$ws.Range("I31").Value = 1824.8462
$ws.Range("J31").Value = 6070.5884
$ws.Range("K31").Value = 1824.8462
$ws.Range("L31").Value = 6070.5884
$ws.Range("M31").Value = -1529.8462
$ws.Range("N31").Value = -6660.5884

$ws.Range("H34").Value = 4230.7666
$ws.Range("I34").Value = 1824.8462
$ws.Range("J34").Value = 6070.5884
$ws.Range("K34").Value = 1824.8462
$ws.Range("L34").Value = 6070.5884
$ws.Range("M34").Value = -1622.8462
$ws.Range("N34").Value = -6474.5884

$ws.Range("H58").Value = 2130.182
$ws.Range("I58").Value = 1711.5172
$ws.Range("J58").Value = 5165.5
$ws.Range("K58").Value = 1711.5172
$ws.Range("L58").Value = 5165.5
$ws.Range("M58").Value = -1508.5172
$ws.Range("N58").Value = -5571.5

$ws.Range("H134").Value = 4494.9375
$ws.Range("I134").Value = 5067.391
$ws.Range("J134").Value = 3032
$ws.Range("K134").Value = 15202.173
$ws.Range("L134").Value = 9096
$ws.Range("M134").Value = -12667.173
$ws.Range("N134").Value = -14166

$ws.Range("H136").Value = 2130.182
$ws.Range("I136").Value = 1711.5172
$ws.Range("J136").Value = 5165.5
$ws.Range("K136").Value = 5134.5516
$ws.Range("L136").Value = 15496.5
$ws.Range("M136").Value = -2584.5516
$ws.Range("N136").Value = -20596.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5674.273
$ws.Range("I56").Value = 5674.273
$ws.Range("K56").Value = 5674.273
$ws.Range("M56").Value = -5144.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35151
$ws.Range("J46").Value = 35151
$ws.Range("L46").Value = 35151
$ws.Range("N46").Value = -35463

$ws.Range("H53").Value = 29999
$ws.Range("J53").Value = 29999
$ws.Range("L53").Value = 29999
$ws.Range("N53").Value = -31261

$ws.Range("H102").Value = 2251.8572
$ws.Range("I102").Value = 1984.7826
$ws.Range("J102").Value = 3480.4
$ws.Range("K102").Value = 1984.7826
$ws.Range("L102").Value = 3480.4
$ws.Range("M102").Value = -362.7826
$ws.Range("N102").Value = -6724.4

$ws.Range("H132").Value = 2865.2144
$ws.Range("I132").Value = 1278
$ws.Range("J132").Value = 4055.625
$ws.Range("K132").Value = 3834
$ws.Range("L132").Value = 12166.875
$ws.Range("M132").Value = -1304
$ws.Range("N132").Value = -17226.875

$ws.Range("H137").Value = 46176
$ws.Range("J137").Value = 46176
$ws.Range("L137").Value = 46176
$ws.Range("N137").Value = -56376

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1195.125
$ws.Range("I82").Value = 797.2857
$ws.Range("J82").Value = 3980
$ws.Range("K82").Value = 797.2857
$ws.Range("L82").Value = 3980
$ws.Range("M82").Value = -436.2857
$ws.Range("N82").Value = -4702

$ws.Range("H85").Value = 1195.125
$ws.Range("I85").Value = 797.2857
$ws.Range("J85").Value = 3980
$ws.Range("K85").Value = 797.2857
$ws.Range("L85").Value = 3980
$ws.Range("M85").Value = 450.7143
$ws.Range("N85").Value = -6476

$ws.Range("H122").Value = 3584.818
$ws.Range("I122").Value = 2491.261
$ws.Range("J122").Value = 6100
$ws.Range("K122").Value = 7473.782999999999
$ws.Range("L122").Value = 18300
$ws.Range("M122").Value = -5023.782999999999
$ws.Range("N122").Value = -23200

$ws.Range("H127").Value = 47800
$ws.Range("J127").Value = 47800
$ws.Range("L127").Value = 47800
$ws.Range("N127").Value = -57720

$ws.Range("H130").Value = 59917.5
$ws.Range("J130").Value = 59917.5
$ws.Range("L130").Value = 59917.5
$ws.Range("N130").Value = -69957.5

$ws.Range("H132").Value = 4166.8447
$ws.Range("I132").Value = 1896.9166
$ws.Range("J132").Value = 7881.273
$ws.Range("K132").Value = 5690.7498
$ws.Range("L132").Value = 23643.819
$ws.Range("M132").Value = -3160.7498
$ws.Range("N132").Value = -28703.819

$ws.Range("H136").Value = 3470.9583
$ws.Range("I136").Value = 1418.5834
$ws.Range("K136").Value = 4255.7502
$ws.Range("M136").Value = -1705.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6806634
$ws.Range("I132").Value = 3869.8572
$ws.Range("J132").Value = 23813544
$ws.Range("K132").Value = 11609.5716
$ws.Range("L132").Value = 71440632
$ws.Range("M132").Value = -9079.571599999999
$ws.Range("N132").Value = -71445692

$ws.Range("H136").Value = 4578.0527
$ws.Range("I136").Value = 1187.1
$ws.Range("J136").Value = 8345.777
$ws.Range("K136").Value = 3561.3
$ws.Range("L136").Value = 25037.331
$ws.Range("M136").Value = -1011.3
$ws.Range("N136").Value = -30137.331
